$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 12: update title (D12) and link (E12)
$ws.Range("D12").Value = "“혼자 공부하는 머신러닝+딥러닝”, “혼자 공부하는 데이터 분석 with 파이썬” 주피터 노트북 재실행 완료"
$ws.Range("E12").Value = "https://tensorflow.blog/2023/11/13/%ed%98%bc%ec%9e%90-%ea%b3%b5%eb%b6%80%ed%95%98%eb%8a%94-%eb%a8%b8%ec%8b%a0%eb%9f%ac%eb%8b%9d%eb%94%a5%eb%9f%ac%eb%8b%9d-%ed%98%bc%ec%9e%90-%ea%b3%b5%eb%b6%80%ed%95%98%eb%8a%94-%eb%8d%b0/"

# Row 28: update title (D28) and link (E28)
$ws.Range("D28").Value = "강화학습 논문 정리 12편 : An Optimistic Perspective on Offline Reinforcement Learning 논문 리뷰 (ICML 2020)"
$ws.Range("E28").Value = "https://ropiens.tistory.com/228"

# Row 51: update title (D51) and link (E51)
$ws.Range("D51").Value = "자바스크립트의 패키지 관리자 npm과 pnpm"
$ws.Range("E51").Value = "https://bskyvision.com/entry/%EC%9E%90%EB%B0%94%EC%8A%A4%ED%81%AC%EB%A6%BD%ED%8A%B8%EC%9D%98-%ED%8C%A8%ED%82%A4%EC%A7%80-%EA%B4%80%EB%A6%AC%EC%9E%90-npm%EA%B3%BC-pnpm"
